$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look like plain numbers (e.g. "1.000",
# "236.25"). Force the column to Text format first so re-assigning these
# values keeps them as literal strings instead of Excel coercing them into
# numeric values (which would drop trailing zeros / reformat the text).
$ws.Range("D2:D51").NumberFormat = "@"

$data = @"
D2|29.167.91
E2|  -0.50%  
D3|1.827.07
E3|  -0.70%  
D4|0.9994
E4|  +0.05%  
D5|236.25
E5|  -1.53%  
D6|0.6026
E6|  -4.02%  
D7|1.000
E7|  +0.03%  
D8|0.07114
E8|  -3.92%  
D9|0.2815
E9|  -2.74%  
D10|24.06
E10|  -2.95%  
D11|0.07649
E11|  -1.07%  
D12|1.831.66
E12|  -0.58%  
D13|4.779
E13|  -3.87%  
D14|0.6407
E14|  -5.41%  
D15|0.000009857
E15|  -2.93%  
B16|WrappedliquidstakedEther2.0
C16|https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth
D16|2.059.08
E16|  -1.64%  
B17|Litecoin
C17|https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc
D17|79.34
E17|  -3.19%  
D18|5.973
E18|  -4.16%  
D19|29.185.75
E19|  -0.50%  
D20|230.99
E20|  +1.09%  
D21|1.000
E21|  +0.08%  
D22|11.73
E22|  -4.59%  
E23|  -4.98%  
D24|1.001
E24|  +0.03%  
D25|155.08
E25|  -2.44%  
D26|8.023
E26|  -5.22%  
D27|0.1280
E27|  -5.34%  
D28|16.65
E28|  -4.44%  
D29|0.06786
E29|  +5.13%  
D30|1.453
E30|  +0.47%  
D31|1.454
E31|  -1.97%  
D32|3.821
E32|  -5.88%  
D33|3.780
E33|  -7.03%  
D34|1.135
E34|  -0.19%  
D35|1.715
E35|  -6.59%  
D36|0.6581
E36|  -5.31%  
D37|2.536
E37|  -1.21%  
D38|1.238.88
E38|  -0.05%  
E39|  -2.03%  
D40|0.01759
E40|  -5.22%  
D41|6.553
E41|  -2.63%  
D42|0.9330
E42|  +0.22%  
D43|0.9996
E43|  +0.06%  
D44|1.972.13
E44|  -2.56%  
D45|100.31
E45|  -0.50%  
D46|63.21
E46|  -3.56%  
E47|  -0.81%  
D48|1.631
E48|  -4.79%  
E49|  -6.77%  
B50|EnergySwap
C50|https://coinranking.com/coin/SbWqqTui-+energyswap-ens
D50|8.516
E50|  -5.10%  
B51|Cronos
C51|https://coinranking.com/coin/65PHZTpmE55b+cronos-cro
D51|0.05581
E51|  -1.59%  
"@

$lines = $data -split "`r`n|`n"
foreach ($line in $lines) {
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|', 2
    $cellRef = $parts[0]
    $value = $parts[1]
    $ws.Range($cellRef).Value = $value
}
